$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reset the scroll/view position (diff removes topLeftCell="A110"; selection stays A2:L121)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Update the rate/derived-value columns (H, J, K, L) for the affected rows
$ws.Range("H25").Value = 0.35
$ws.Range("J25").Value = 17.260273972602739
$ws.Range("K25").Value = 124.88219178082193
$ws.Range("L25").Value = 0.50455555555555553
$ws.Range("H26").Value = 0.38
$ws.Range("J26").Value = 28.109589041095891
$ws.Range("K26").Value = 137.71780821917807
$ws.Range("L26").Value = 0.50988888888888884
$ws.Range("H27").Value = 0.38
$ws.Range("J27").Value = 37.479452054794521
$ws.Range("K27").Value = 151.22739726027399
$ws.Range("L27").Value = 0.51938888888888901
$ws.Range("H28").Value = 0.36
$ws.Range("J28").Value = 17.753424657534246
$ws.Range("K28").Value = 125.40739726027397
$ws.Range("L28").Value = 0.51520555555555536
$ws.Range("H29").Value = 0.39
$ws.Range("J29").Value = 28.849315068493155
$ws.Range("K29").Value = 138.51301369863015
$ws.Range("L29").Value = 0.52063888888888887
$ws.Range("H30").Value = 0.39
$ws.Range("J30").Value = 38.465753424657535
$ws.Range("K30").Value = 152.31232876712329
$ws.Range("L30").Value = 0.53038888888888891
$ws.Range("H44").Value = 0.44
$ws.Range("J44").Value = 21.698630136986303
$ws.Range("K44").Value = 122.30712328767122
$ws.Range("L44").Value = 0.45233888888888857
$ws.Range("H45").Value = 0.44
$ws.Range("J45").Value = 32.547945205479458
$ws.Range("K45").Value = 133.21068493150685
$ws.Range("L45").Value = 0.44895925925925922
$ws.Range("H46").Value = 0.44
$ws.Range("J46").Value = 43.397260273972606
$ws.Range("K46").Value = 144.11424657534243
$ws.Range("L46").Value = 0.44726944444444422
$ws.Range("H49").Value = 0.38
$ws.Range("J49").Value = 37.479452054794521
$ws.Range("K49").Value = 144.35342465753425
$ws.Range("L49").Value = 0.44969444444444445
$ws.Range("H50").Value = 0.21
$ws.Range("J50").Value = 10.356164383561643
$ws.Range("K50").Value = 121.39178082191782
$ws.Range("L50").Value = 0.43377777777777793
$ws.Range("H51").Value = 0.27
$ws.Range("J51").Value = 19.972602739726032
$ws.Range("K51").Value = 131.96986301369864
$ws.Range("L51").Value = 0.43218518518518528
$ws.Range("H52").Value = 0.31
$ws.Range("J52").Value = 30.575342465753426
$ws.Range("K52").Value = 143.63287671232877
$ws.Range("L52").Value = 0.44238888888888889
$ws.Range("H107").Value = 0.52
$ws.Range("J107").Value = 8.5479452054794507
$ws.Range("K107").Value = 110.71890410958905
$ws.Range("L107").Value = 0.65206666666666691
$ws.Range("H108").Value = 0.52
$ws.Range("J108").Value = 12.821917808219178
$ws.Range("K108").Value = 115.07835616438356
$ws.Range("L108").Value = 0.61151111111111078
$ws.Range("H109").Value = 0.52
$ws.Range("J109").Value = 25.643835616438356
$ws.Range("K109").Value = 128.15671232876713
$ws.Range("L109").Value = 0.57095555555555566
$ws.Range("H110").Value = 0.52
$ws.Range("J110").Value = 38.465753424657535
$ws.Range("K110").Value = 141.23506849315066
$ws.Range("L110").Value = 0.55743703703703673
$ws.Range("H111").Value = 0.52
$ws.Range("J111").Value = 51.287671232876711
$ws.Range("K111").Value = 154.31342465753423
$ws.Range("L111").Value = 0.55067777777777771

Write-Output "Updated H/J/K/L values for rows 25-30, 44-46, 49-52, 107-111"
